$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original text formatting so that
# numeric-looking values (e.g. "5.45", "0.0993") are not auto-converted
# into floating point numbers by Excels input parser.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '59.419.75'
$ws.Range("E2").Value = '  -5.64%  '
$ws.Range("D3").Value = '2.447.07'
$ws.Range("E3").Value = '  -8.77%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '535.82'
$ws.Range("E5").Value = '  -2.53%  '
$ws.Range("E6").Value = '  -6.28%  '
$ws.Range("E7").Value = '  -0.23%  '
$ws.Range("E8").Value = '  -3.26%  '
$ws.Range("D9").Value = '0.0993'
$ws.Range("E9").Value = '  -5.81%  '
$ws.Range("E10").Value = '  -2.29%  '
$ws.Range("D11").Value = '5.45'
$ws.Range("E11").Value = '  +6.87%  '
$ws.Range("D12").Value = '0.351'
$ws.Range("E12").Value = '  -4.48%  '
$ws.Range("D13").Value = '2.881.92'
$ws.Range("E13").Value = '  -8.73%  '
$ws.Range("D14").Value = '24.17'
$ws.Range("E14").Value = '  -7.18%  '
$ws.Range("D15").Value = '59.390.17'
$ws.Range("E15").Value = '  -5.48%  '
$ws.Range("E16").Value = '  -5.76%  '
$ws.Range("D17").Value = '2.494.40'
$ws.Range("E17").Value = '  -7.09%  '
$ws.Range("D18").Value = '11.20'
$ws.Range("E18").Value = '  -5.75%  '
$ws.Range("D19").Value = '4.37'
$ws.Range("E19").Value = '  -4.35%  '
$ws.Range("D20").Value = '325.06'
$ws.Range("E20").Value = '  -5.22%  '
$ws.Range("D21").Value = '0.972'
$ws.Range("E21").Value = '  -2.61%  '
$ws.Range("D22").Value = '5.72'
$ws.Range("E22").Value = '  -9.53%  '
$ws.Range("D23").Value = '0.468'
$ws.Range("E23").Value = '  -7.05%  '
$ws.Range("D24").Value = '60.42'
$ws.Range("E24").Value = '  -4.75%  '
$ws.Range("E25").Value = '  -3.68%  '
$ws.Range("E26").Value = '  -2.45%  '
$ws.Range("E27").Value = '  -4.62%  '
$ws.Range("D28").Value = '1.31'
$ws.Range("E28").Value = '  -1.97%  '
$ws.Range("E29").Value = '  -1.86%  '
$ws.Range("E30").Value = '  -5.23%  '
$ws.Range("E31").Value = '  -9.97%  '
$ws.Range("D32").Value = '0.998'
$ws.Range("E32").Value = '  -0.08%  '
$ws.Range("D33").Value = '156.02'
$ws.Range("E33").Value = '  -6.64%  '
$ws.Range("D34").Value = '4.55'
$ws.Range("E34").Value = '  -5.19%  '
$ws.Range("D35").Value = '18.38'
$ws.Range("E35").Value = '  -5.76%  '
$ws.Range("E36").Value = '  -5.05%  '
$ws.Range("D37").Value = '1.75'
$ws.Range("E37").Value = '  -0.70%  '
$ws.Range("D38").Value = '314.81'
$ws.Range("E38").Value = '  -6.90%  '
$ws.Range("E39").Value = '  -6.49%  '
$ws.Range("E40").Value = '  -8.19%  '
$ws.Range("D41").Value = '36.88'
$ws.Range("E41").Value = '  -3.61%  '
$ws.Range("D42").Value = '3.72'
$ws.Range("E42").Value = '  -5.51%  '
$ws.Range("D43").Value = '0.996'
$ws.Range("E43").Value = '  -0.28%  '
$ws.Range("D44").Value = '10.72'
$ws.Range("E44").Value = '  -2.92%  '
$ws.Range("D45").Value = '0.583'
$ws.Range("E45").Value = '  -5.39%  '
$ws.Range("D46").Value = '0.0940'
$ws.Range("E46").Value = '  -3.20%  '
$ws.Range("D47").Value = '0.0526'
$ws.Range("E47").Value = '  -6.16%  '
$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").Value = '19.09'
$ws.Range("E48").Value = '  -7.68%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '18.66'
$ws.Range("E49").Value = '  -7.97%  '
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").Value = '0.0230'
$ws.Range("E50").Value = '  -3.93%  '
$ws.Range("D51").Value = '1.989.56'
$ws.Range("E51").Value = '  -4.72%  '
